$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (bold, bordered, centered) from H1
# onto the new header cells I1:J1 so the new headers match the look
# of the existing ones (B1:H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I (I0) and J (IF)
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 7

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 2

$ws.Range("I5").Value = 9
$ws.Range("J5").Value = 9

$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 5

$ws.Range("I7").Value = 4
$ws.Range("J7").Value = 5
